# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# D-column "Price" and E-column "Volume(1h)" cells are plain text in the
# source sheet (t="inlineStr"), so purely numeric-looking Price values are
# set with a leading apostrophe to force Excel to keep them as text instead
# of auto-coercing to a Number (this mirrors typing '212.62 into the cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.558.31'
$ws.Range("E2").Value = '  +2.96%  '
$ws.Range("D3").Value = '1.605.36'
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'212.62"
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("E6").Value = '  +7.25%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'26.93"
$ws.Range("E8").Value = '  +7.32%  '
$ws.Range("D9").Value = "'43.46"
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  +2.45%  '
$ws.Range("E11").Value = '  +2.75%  '
$ws.Range("D13").Value = '1.833.89'
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").Value = '1.587.62'
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("E15").Value = '  +3.80%  '
$ws.Range("D16").Value = '29.558.15'
$ws.Range("E16").Value = '  +3.02%  '
$ws.Range("E18").Value = '  +3.40%  '
$ws.Range("D19").Value = "'240.43"
$ws.Range("E19").Value = '  +4.56%  '
$ws.Range("E20").Value = '  +3.41%  '
$ws.Range("E21").Value = '  +1.83%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("E24").Value = '  +2.25%  '
$ws.Range("E25").Value = '  +0.69%  '
$ws.Range("D26").Value = "'154.55"
$ws.Range("E26").Value = '  +2.00%  '
$ws.Range("E27").Value = '  +5.17%  '
$ws.Range("D28").Value = "'15.26"
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("D29").Value = "'6.40"
$ws.Range("E29").Value = '  +2.55%  '
$ws.Range("D30").Value = "'0.998"
$ws.Range("E31").Value = '  +2.72%  '
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("E34").Value = '  +3.44%  '
$ws.Range("D35").Value = '1.408.81'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("E37").Value = '  +4.32%  '
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("E40").Value = '  +2.51%  '
$ws.Range("D41").Value = "'0.538"
$ws.Range("E41").Value = '  +3.80%  '
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("E43").Value = '  +7.20%  '
$ws.Range("D44").Value = "'53.61"
$ws.Range("E44").Value = '  +24.55%  '
$ws.Range("E45").Value = '  +3.30%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = "'65.82"
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").Value = '1.744.74'
$ws.Range("E49").Value = '  +2.61%  '
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").Value = "'86.67"
$ws.Range("E51").Value = '  +1.94%  '
